$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# Column A: blank note field (source row has an empty value but the cell
# itself is present). Force the cell to materialize even though its
# content is empty.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = ""
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "احمد"

# Column C: quantity looks numeric ("22") but must stay text, matching the
# source data where every column is stored as text.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "22"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"
$ws.Cells.Item($row, 6).Value = "C3"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٠٣:٥٢ م"
